$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    "B" = 0.9999845286516351
    "C" = 0.9991389381953254
    "D" = 0.9999955598836291
    "E" = 0.9999999999992062
    "F" = 0.9999965254174424
    "G" = 0.00001444182458493536
    "H" = 0.000803763398418822
    "I" = 0.000004467702705704711
    "J" = 0.0000000000001175802112053254
    "K" = 0.000002233851392423344
    "L" = 0.0001971186733916371
    "M" = 0.003800240069381849
    "N" = 1.00001125188972
    "O" = 0.003962024156716867
    "P" = 136.2907641523688
    "Q" = 205.7666861698563
}

for ($row = 2; $row -le 26; $row++) {
    foreach ($col in $newValues.Keys) {
        $ws.Range("$col$row").Value = $newValues[$col]
    }
}

$wb.Save()
